# Apply weekly update to the "Hortaliza, Mapocho Venta Directa de Santiago - Esparragos" sheet.
# The edit inserts three new price records (rows 5-7) ahead of the existing data,
# pushing the previously-existing rows 5-50 down to rows 8-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 5. This shifts old rows 5..50 down to 8..53,
# and the new blank rows inherit the date-formatted style from row 4 (column D),
# matching the existing pattern used throughout the sheet.
$ws.Range("A5:A7").EntireRow.Insert()

function Set-Row($r, $d, $h, $i, $j, $k, $l, $m, $n, $o, $p, $q) {
    $ws.Cells.Item($r, 1).Value = 12
    $ws.Cells.Item($r, 2).Value = "Mapocho Venta Directa de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = 300000000
    $ws.Cells.Item($r, 7).Value = "Espárragos"
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}

# New row 5: Banquete, 2021-10-21 (serial 44490)
Set-Row 5 44490 "Sin especificar" "Banquete" 300 1200 1200 1200 "`$/kilo" "Región Metropolitana" 1200 1

# New row 6: Primera, 2021-10-21 (serial 44490)
Set-Row 6 44490 "Sin especificar" "Primera" 330 1000 1000 1000 "`$/kilo" "Región Metropolitana" 1000 1

# New row 7: Segunda, 2021-10-21 (serial 44490)
Set-Row 7 44490 "Sin especificar" "Segunda" 350 800 800 800 "`$/kilo" "Región Metropolitana" 800 1
